$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A269").Value = 'FSD'
$ws.Range("B269").Value = 'Sioux Falls, South Dakota'
$ws.Range("C269").Value = 43.540819819502
$ws.Range("D269").Value = -96.65511577730963
$ws.Range("E269").Value = 'US'
$ws.Range("F269").Value = 'North America'
$ws.Range("G269").Value = 'Sioux Falls'

$ws.Range("A270").Value = 'STL'
$ws.Range("B270").Value = 'St. Louis, MO, United States'
$ws.Range("C270").Value = 38.7486991882
$ws.Range("D270").Value = -90.37000274659999
$ws.Range("E270").Value = 'US'
$ws.Range("F270").Value = 'North America'
$ws.Range("G270").Value = 'St. Louis'

$ws.Range("A271").Value = 'TLH'
$ws.Range("B271").Value = 'Tallahassee, FL, United States'
$ws.Range("C271").Value = 30.3964996338
$ws.Range("D271").Value = -84.3503036499
$ws.Range("E271").Value = 'US'
$ws.Range("F271").Value = 'North America'
$ws.Range("G271").Value = 'Tallahassee'

$ws.Range("A272").Value = 'TPA'
$ws.Range("B272").Value = 'Tampa, FL, United States'
$ws.Range("C272").Value = 27.9755001068
$ws.Range("D272").Value = -82.533203125
$ws.Range("E272").Value = 'US'
$ws.Range("F272").Value = 'North America'
$ws.Range("G272").Value = 'Tampa'

$ws.Range("A273").Value = 'YYZ'
$ws.Range("B273").Value = 'Toronto, ON, Canada'
$ws.Range("C273").Value = 43.6772003174
$ws.Range("D273").Value = -79.63059997560001
$ws.Range("E273").Value = 'CA'
$ws.Range("F273").Value = 'North America'
$ws.Range("G273").Value = 'Toronto'

$ws.Range("A274").Value = 'YVR'
$ws.Range("B274").Value = 'Vancouver, BC, Canada'
$ws.Range("C274").Value = 49.193901062
$ws.Range("D274").Value = -123.183998108
$ws.Range("E274").Value = 'CA'
$ws.Range("F274").Value = 'North America'
$ws.Range("G274").Value = 'Vancouver'

$ws.Range("A275").Value = 'YWG'
$ws.Range("B275").Value = 'Winnipeg, MB, Canada'
$ws.Range("C275").Value = 49.9099998474
$ws.Range("D275").Value = -97.2398986816
$ws.Range("E275").Value = 'CA'
$ws.Range("F275").Value = 'North America'
$ws.Range("G275").Value = 'Winnipeg'

$ws.Range("A276").Value = 'ADL'
$ws.Range("B276").Value = 'Adelaide, SA, Australia'
$ws.Range("C276").Value = -34.9431729
$ws.Range("D276").Value = 138.5335637
$ws.Range("E276").Value = 'AU'
$ws.Range("F276").Value = 'Oceania'
$ws.Range("G276").Value = 'Adelaide'

$ws.Range("A277").Value = 'AKL'
$ws.Range("B277").Value = 'Auckland, New Zealand'
$ws.Range("C277").Value = -37.0080986023
$ws.Range("D277").Value = 174.792007446
$ws.Range("E277").Value = 'NZ'
$ws.Range("F277").Value = 'Oceania'
$ws.Range("G277").Value = 'Auckland'

$ws.Range("A278").Value = 'BNE'
$ws.Range("B278").Value = 'Brisbane, QLD, Australia'
$ws.Range("C278").Value = -27.3841991425
$ws.Range("D278").Value = 153.117004394
$ws.Range("E278").Value = 'AU'
$ws.Range("F278").Value = 'Oceania'
$ws.Range("G278").Value = 'Brisbane'

$ws.Range("A279").Value = 'CBR'
$ws.Range("B279").Value = 'Canberra, ACT, Australia'
$ws.Range("C279").Value = -35.3069000244
$ws.Range("D279").Value = 149.1950073242
$ws.Range("E279").Value = 'AU'
$ws.Range("F279").Value = 'Oceania'
$ws.Range("G279").Value = 'Canberra'

$ws.Range("A280").Value = 'CHC'
$ws.Range("B280").Value = 'Christchurch, New Zealand'
$ws.Range("C280").Value = -43.4893989563
$ws.Range("D280").Value = 172.5319976807
$ws.Range("E280").Value = 'NZ'
$ws.Range("F280").Value = 'Oceania'
$ws.Range("G280").Value = 'Christchurch'

$ws.Range("A281").Value = 'GUM'
$ws.Range("B281").Value = 'Hagatna, Guam'
$ws.Range("C281").Value = 13.4834003448
$ws.Range("D281").Value = 144.796005249
$ws.Range("E281").Value = 'GU'
$ws.Range("F281").Value = 'Asia Pacific'
$ws.Range("G281").Value = 'Hagatna'

$ws.Range("A282").Value = 'HBA'
$ws.Range("B282").Value = 'Hobart, Australia'
$ws.Range("C282").Value = -42.883209
$ws.Range("D282").Value = 147.331665
$ws.Range("E282").Value = 'AU'
$ws.Range("F282").Value = 'Oceania'
$ws.Range("G282").Value = 'Hobart'

$ws.Range("A283").Value = 'MEL'
$ws.Range("B283").Value = 'Melbourne, VIC, Australia'
$ws.Range("C283").Value = -37.6733016968
$ws.Range("D283").Value = 144.843002319
$ws.Range("E283").Value = 'AU'
$ws.Range("F283").Value = 'Oceania'
$ws.Range("G283").Value = 'Melbourne'

$ws.Range("A284").Value = 'NOU'
$ws.Range("B284").Value = 'Noumea, New Caledonia'
$ws.Range("C284").Value = -22.0146007538
$ws.Range("D284").Value = 166.212997436
$ws.Range("E284").Value = 'NC'
$ws.Range("F284").Value = 'Oceania'
$ws.Range("G284").Value = 'Noumea'

$ws.Range("A285").Value = 'PER'
$ws.Range("B285").Value = 'Perth, WA, Australia'
$ws.Range("C285").Value = -31.9402999878
$ws.Range("D285").Value = 115.967002869
$ws.Range("E285").Value = 'AU'
$ws.Range("F285").Value = 'Oceania'
$ws.Range("G285").Value = 'Perth'

$ws.Range("A286").Value = 'SYD'
$ws.Range("B286").Value = 'Sydney, NSW, Australia'
$ws.Range("C286").Value = -33.9460983276
$ws.Range("D286").Value = 151.177001953
$ws.Range("E286").Value = 'AU'
$ws.Range("F286").Value = 'Oceania'
$ws.Range("G286").Value = 'Sydney'

$ws.Range("A287").Value = 'PPT'
$ws.Range("B287").Value = 'Tahiti, French Polynesia'
$ws.Range("C287").Value = -17.5536994934
$ws.Range("D287").Value = -149.606994629
$ws.Range("E287").Value = 'PF'
$ws.Range("F287").Value = 'Oceania'
$ws.Range("G287").Value = 'Tahiti'

$ws.Range("A288").Value = 'REC'
$ws.Range("B288").Value = 'Recife, Brazil'
$ws.Range("C288").Value = -8.126489639300001
$ws.Range("D288").Value = -34.9235992432
$ws.Range("E288").Value = 'BR'
$ws.Range("F288").Value = 'South America'
$ws.Range("G288").Value = 'Recife'

$ws.Range("A289").Value = 'STR'
$ws.Range("B289").Value = 'Stuttgart, Germany'
$ws.Range("C289").Value = 48.783333
$ws.Range("D289").Value = 9.183332999999999
$ws.Range("E289").Value = 'DE'
$ws.Range("F289").Value = 'Europe'
$ws.Range("G289").Value = 'Stuttgart'

$ws.Range("A290").Value = 'COK'
$ws.Range("B290").Value = 'Kochi, India'
$ws.Range("C290").Value = 9.9312
$ws.Range("D290").Value = 76.26730000000001
$ws.Range("E290").Value = 'IN'
$ws.Range("F290").Value = 'Asia Pacific'
$ws.Range("G290").Value = 'Kochi'

$ws.Range("A291").Value = 'FUK'
$ws.Range("B291").Value = 'Fukuoka, Japan'
$ws.Range("C291").Value = 33.5902
$ws.Range("D291").Value = 130.4017
$ws.Range("E291").Value = 'JP'
$ws.Range("F291").Value = 'Asia Pacific'
$ws.Range("G291").Value = 'Fukuoka'

$ws.Range("A292").Value = 'ADB'
$ws.Range("B292").Value = 'Izmir, Turkey'
$ws.Range("C292").Value = 38.32377
$ws.Range("D292").Value = 27.14317
$ws.Range("E292").Value = 'TR'
$ws.Range("F292").Value = 'Europe'
$ws.Range("G292").Value = 'Izmir'

$ws.Range("A293").Value = 'SFO'
$ws.Range("B293").Value = 'San Francisco, United States'
$ws.Range("C293").Value = 37.6189994812
$ws.Range("D293").Value = -122.375
$ws.Range("E293").Value = 'US'
$ws.Range("F293").Value = 'North America'
$ws.Range("G293").Value = 'San Francisco'
